$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-18 changes from serial date 45174 (2023-09-05)
# to 45175 (2023-09-06).
foreach ($row in 2..18) {
    $ws.Cells.Item($row, 3).Value = 45175
}
